# Sales & Redemption report: rename sheet, refresh reporting period,
# restructure the summary header block (B/C columns -> A/B columns),
# and tweak a couple of column widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Rename the worksheet tab
# ------------------------------------------------------------------
$ws.Name = "Sales"

# ------------------------------------------------------------------
# 2. Shift the report-summary block one column to the left:
#    B1          -> A1
#    B2:C2       -> A2:B2
#    B3:C3       -> A3:B3
#    B4:C4       -> A4:B4   (also refresh the reporting-period text)
#    B5:C5       -> A5:B5   (amount cell becomes left aligned)
#    B6:C6       -> A6:B6   (count cell becomes left aligned, plain number fmt)
# ------------------------------------------------------------------

# Row 1 - title
$ws.Range("B1").Copy($ws.Range("A1"))
$ws.Range("B1").Clear()

# Row 2 - Store name
$ws.Range("B2").Copy($ws.Range("A2"))
$ws.Range("C2").Copy($ws.Range("B2"))
$ws.Range("C2").Clear()

# Row 3 - Product name
$ws.Range("B3").Copy($ws.Range("A3"))
$ws.Range("C3").Copy($ws.Range("B3"))
$ws.Range("C3").Clear()

# Row 4 - Reporting period (also update the end date 22-Jun-2024 -> 24-Jun-2024)
$ws.Range("B4").Copy($ws.Range("A4"))
$ws.Range("C4").Copy($ws.Range("B4"))
$ws.Range("C4").Clear()
$ws.Range("B4").Value = "01-Jul-2023 To 24-Jun-2024"

# Row 5 - GGC total net amount
$ws.Range("B5").Copy($ws.Range("A5"))
$ws.Range("C5").Copy($ws.Range("B5"))
$ws.Range("C5").Clear()
$ws.Range("B5").HorizontalAlignment = -4131   # xlLeft

# Row 6 - Count
$ws.Range("B6").Copy($ws.Range("A6"))
$ws.Range("C6").Copy($ws.Range("B6"))
$ws.Range("C6").Clear()
$ws.Range("B6").HorizontalAlignment = -4131   # xlLeft

# ------------------------------------------------------------------
# 3. Column width tweaks (column A wider, column C narrower)
# ------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 25.96
$ws.Columns.Item(3).ColumnWidth = 19.53

Write-Host "Sales reporting workbook updated"
